# Lesson 6.5 Lists of Lists — apply commit:
# "changed colors in L6.3 for the color-blind. Put in refs to files throughout"
#
# Concretely this touches:
#  1) The cached Header/Footer date text (10/17/2015 -> 10/19/2015) that is
#     stamped on the slide master, every slide layout, and the notes master.
#  2) Slide 22: merge "Sexp"+"OfX" -> "SexpOfX" and
#     "ListOfSexp"+"Of"+"X" -> "ListOfSexpOfX" (same run formatting, so the
#     runs collapse into one run each).
#  3) Slide 28: add a new first bullet pointing at the example file, and
#     split the "If you have questions..." sentence so "If " is its own run.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Update the cached date field text everywhere it is stamped.
# ---------------------------------------------------------------------
function Set-DateShapeText($shapes, $newText) {
  for ($i = 1; $i -le $shapes.Count; $i++) {
    $sh = $shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
      $sh.TextFrame.TextRange.Text = $newText
    }
  }
}

# Slide master date placeholder
Set-DateShapeText $p.SlideMaster.Shapes "10/19/2015"

# Every slide layout's date placeholder
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
  Set-DateShapeText $layouts.Item($li).Shapes "10/19/2015"
}

# Notes master date placeholder is an auto-updating field; it only accepts
# edits through the HeadersFooters API.
$p.NotesMaster.HeadersFooters.DateAndTime.Text = "10/19/2015"

# ---------------------------------------------------------------------
# 2) Slide 22 ("The Template for SexpX"): merge split identifier runs.
# ---------------------------------------------------------------------
$s22 = $p.Slides.Item(22)
$sh22 = $s22.Shapes.Item(2)
$tf22 = $sh22.TextFrame

$trSexp = $tf22.TextRange.Find("SexpOfX", 0)
$trSexp.Text = "SexpOfX"

$trListSexp = $tf22.TextRange.Find("ListOfSexpOfX", 0)
$trListSexp.Text = "ListOfSexpOfX"

# ---------------------------------------------------------------------
# 3) Slide 28 ("Next Steps"): add the file reference bullet and split the
#    discussion-board sentence.
# ---------------------------------------------------------------------
$s28 = $p.Slides.Item(28)
$sh28 = $s28.Shapes.Item(2)
$tf28 = $sh28.TextFrame

$firstPara = $tf28.TextRange.Paragraphs(1, 1)
$firstPara.InsertBefore("Study the file 06-5-sos-and-loss.rkt in the Examples folder`r") | Out-Null

$examplesRange = $tf28.TextRange.Find("Examples folder", 0)
$examplesRange.Text = "Examples folder"

$ifRange = $tf28.TextRange.Find("If ", 0)
$ifRange.Text = "If "
